# LOM3116.xlsx — update teaching staff / activation-date cells.
#
# Net effect (per the target diff):
#   B10/C10  "Aulas expositivas teóricas..."        -> "471420 - Carlos Antonio Reis Pereira Baptista"
#   B13/C13  "90 h"                                  -> "01/01/2023"
#   B15/C15  "Semestral"                              -> "471420 - Carlos Antonio Reis Pereira Baptista"
#   B18/C18  "01/01/2023"                             -> "3586455 - Cassius Olivio Figueiredo Terra Ruchert"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: Docentes responsáveis (Objetivos block) ---
$ws.Range("B10").Value = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Range("C10").Value = "471420 - Carlos Antonio Reis Pereira Baptista"

# --- Row 13: Ativação date. Copy from B8/C8, which already hold the
# literal text "01/01/2023", so the cell keeps its text type/style
# instead of Excel auto-converting the string into a date serial. ---
$ws.Range("B8").Copy($ws.Range("B13"))
$ws.Range("C8").Copy($ws.Range("C13"))

# --- Row 15: Programa (second teacher credit) ---
$ws.Range("B15").Value = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Range("C15").Value = "471420 - Carlos Antonio Reis Pereira Baptista"

# --- Row 18: Método (responsible professor) ---
$ws.Range("B18").Value = "3586455 - Cassius Olivio Figueiredo Terra Ruchert"
$ws.Range("C18").Value = "3586455 - Cassius Olivio Figueiredo Terra Ruchert"
